$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3000.348
$ws.Range("J137").Value = 3799.5833
$ws.Range("L137").Value = 11398.7499
$ws.Range("N137").Value = -16498.7499
$ws.Range("H138").Value = 2778.1128
$ws.Range("I138").Value = 1890.3846
$ws.Range("J138").Value = 3419.25
$ws.Range("K138").Value = 5671.1538
$ws.Range("L138").Value = 10257.75
$ws.Range("M138").Value = -531.1538
$ws.Range("N138").Value = -20537.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 10955.182
$ws.Range("J2").Value = 13090
$ws.Range("L2").Value = 13090
$ws.Range("N2").Value = -13316
$ws.Range("H21").Value = 578.3333
$ws.Range("I21").Value = 562.5714
$ws.Range("J21").Value = 799
$ws.Range("K21").Value = 562.5714
$ws.Range("L21").Value = 799
$ws.Range("M21").Value = -188.5714
$ws.Range("N21").Value = -1547
$ws.Range("H30").Value = 4878.2
$ws.Range("I30").Value = 1095.25
$ws.Range("J30").Value = 20010
$ws.Range("K30").Value = 1095.25
$ws.Range("L30").Value = 20010
$ws.Range("M30").Value = -945.25
$ws.Range("N30").Value = -20310
$ws.Range("H32").Value = 3602.932
$ws.Range("I32").Value = 3059.1025
$ws.Range("J32").Value = 7844.8
$ws.Range("K32").Value = 3059.1025
$ws.Range("L32").Value = 7844.8
$ws.Range("M32").Value = -2772.1025
$ws.Range("N32").Value = -8418.799999999999
$ws.Range("H37").Value = 20018.5
$ws.Range("J37").Value = 30038
$ws.Range("L37").Value = 30038
$ws.Range("N37").Value = -30584
$ws.Range("H61").Value = 6851.5776
$ws.Range("I61").Value = 5349.5
$ws.Range("J61").Value = 11494.363
$ws.Range("K61").Value = 5349.5
$ws.Range("L61").Value = 11494.363
$ws.Range("M61").Value = -5137.5
$ws.Range("N61").Value = -11918.363
$ws.Range("H116").Value = 10955.182
$ws.Range("J116").Value = 13090
$ws.Range("L116").Value = 13090
$ws.Range("N116").Value = -17678
$ws.Range("H136").Value = 6851.5776
$ws.Range("I136").Value = 5349.5
$ws.Range("J136").Value = 11494.363
$ws.Range("K136").Value = 16048.5
$ws.Range("L136").Value = 34483.089
$ws.Range("M136").Value = -13498.5
$ws.Range("N136").Value = -39583.089
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H139").Value = 70532.336
$ws.Range("J139").Value = 70532.336
$ws.Range("L139").Value = 70532.336
$ws.Range("N139").Value = -80812.336
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("N140").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 10955.182
$ws.Range("J3").Value = 13090
$ws.Range("L3").Value = 13090
$ws.Range("N3").Value = -13318
$ws.Range("H86").Value = 4384.273
$ws.Range("I86").Value = 3024.5557
$ws.Range("J86").Value = 10503
$ws.Range("K86").Value = 3024.5557
$ws.Range("L86").Value = 10503
$ws.Range("M86").Value = -1901.5557
$ws.Range("N86").Value = -12749
$ws.Range("H89").Value = 4384.273
$ws.Range("I89").Value = 3024.5557
$ws.Range("J89").Value = 10503
$ws.Range("K89").Value = 15122.7785
$ws.Range("L89").Value = 52515
$ws.Range("M89").Value = -9506.7785
$ws.Range("N89").Value = -63747
$ws.Range("H105").Value = 13277.174
$ws.Range("I105").Value = 16357.143
$ws.Range("K105").Value = 16357.143
$ws.Range("M105").Value = -14610.143
$ws.Range("H107").Value = 199
$ws.Range("I107").Value = 199
$ws.Range("K107").Value = 199
$ws.Range("M107").Value = 1721
$ws.Range("H134").Value = 2643.5527
$ws.Range("J134").Value = 11655.167
$ws.Range("L134").Value = 34965.501
$ws.Range("N134").Value = -40035.501

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5346.423
$ws.Range("I58").Value = 3865.2354
$ws.Range("J58").Value = 8144.222
$ws.Range("K58").Value = 3865.2354
$ws.Range("L58").Value = 8144.222
$ws.Range("M58").Value = -3662.2354
$ws.Range("N58").Value = -8550.222
$ws.Range("H86").Value = 14220.777
$ws.Range("J86").Value = 15497.833
$ws.Range("L86").Value = 15497.833
$ws.Range("N86").Value = -17743.833
$ws.Range("H89").Value = 14220.777
$ws.Range("J89").Value = 15497.833
$ws.Range("L89").Value = 77489.16500000001
$ws.Range("N89").Value = -88721.16500000001
$ws.Range("H103").Value = 41428.715
$ws.Range("I103").Value = 35000.25
$ws.Range("J103").Value = 50000
$ws.Range("K103").Value = 35000.25
$ws.Range("L103").Value = 50000
$ws.Range("M103").Value = -33828.25
$ws.Range("N103").Value = -52344
$ws.Range("H132").Value = 3392.3696
$ws.Range("I132").Value = 2965.7297
$ws.Range("K132").Value = 8897.1891
$ws.Range("M132").Value = -6367.1891
$ws.Range("H134").Value = 2669.6316
$ws.Range("I134").Value = 1633.3572
$ws.Range("J134").Value = 5571.2
$ws.Range("K134").Value = 4900.071599999999
$ws.Range("L134").Value = 16713.6
$ws.Range("M134").Value = -2365.071599999999
$ws.Range("N134").Value = -21783.6
$ws.Range("H136").Value = 5346.423
$ws.Range("I136").Value = 3865.2354
$ws.Range("J136").Value = 8144.222
$ws.Range("K136").Value = 11595.7062
$ws.Range("L136").Value = 24432.666
$ws.Range("M136").Value = -9045.706200000001
$ws.Range("N136").Value = -29532.666
$ws.Range("H137").Value = 70000
$ws.Range("J137").Value = 70000
$ws.Range("L137").Value = 70000
$ws.Range("N137").Value = -80200
$ws.Range("H138").Value = 54997.2
$ws.Range("J138").Value = 54997.2
$ws.Range("L138").Value = 54997.2
$ws.Range("N138").Value = -65277.2
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H140").Value = 69999
$ws.Range("J140").Value = 69999
$ws.Range("L140").Value = 69999
$ws.Range("N140").Value = -80359
$ws.Range("H141").Value = 237193
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 237193
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -247553

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 30158.35
$ws.Range("I2").Value = 146.57143
$ws.Range("J2").Value = 100185.836
$ws.Range("K2").Value = 879.42858
$ws.Range("L2").Value = 601115.0159999999
$ws.Range("M2").Value = -766.42858
$ws.Range("N2").Value = -601341.0159999999
$ws.Range("H80").Value = 6001.75
$ws.Range("I80").Value = 3502
$ws.Range("K80").Value = 10506
$ws.Range("M80").Value = -9570
$ws.Range("H81").Value = 146792
$ws.Range("J81").Value = 339666.66
$ws.Range("L81").Value = 1018999.98
$ws.Range("N81").Value = -1021245.98
$ws.Range("H83").Value = 6001.75
$ws.Range("I83").Value = 3502
$ws.Range("K83").Value = 31518
$ws.Range("M83").Value = -26838
$ws.Range("H84").Value = 146792
$ws.Range("J84").Value = 339666.66
$ws.Range("L84").Value = 3056999.94
$ws.Range("N84").Value = -3068231.94
$ws.Range("H107").Value = 1377.6471
$ws.Range("I107").Value = 1960
$ws.Range("J107").Value = 722.5
$ws.Range("K107").Value = 5880
$ws.Range("L107").Value = 2167.5
$ws.Range("M107").Value = -3960
$ws.Range("N107").Value = -6007.5
$ws.Range("H113").Value = 1261.2222
$ws.Range("I113").Value = 992.8333
$ws.Range("J113").Value = 1798
$ws.Range("K113").Value = 2978.4999
$ws.Range("L113").Value = 5394
$ws.Range("M113").Value = -808.4998999999998
$ws.Range("N113").Value = -9734
$ws.Range("H131").Value = 7779507.5
$ws.Range("I131").Value = 1924.4
$ws.Range("J131").Value = 12964563
$ws.Range("K131").Value = 5773.200000000001
$ws.Range("L131").Value = 38893689
$ws.Range("M131").Value = -733.2000000000007
$ws.Range("N131").Value = -38903769
$ws.Range("H134").Value = 2332
$ws.Range("I134").Value = 1920.8
$ws.Range("J134").Value = 8500
$ws.Range("K134").Value = 5762.4
$ws.Range("L134").Value = 25500
$ws.Range("M134").Value = -692.3999999999996
$ws.Range("N134").Value = -35640

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 51286
$ws.Range("I132").Value = 60376.61
$ws.Range("K132").Value = 181129.83
$ws.Range("M132").Value = -178599.83
$ws.Range("H135").Value = 59949
$ws.Range("J135").Value = 59949
$ws.Range("L135").Value = 59949
$ws.Range("N135").Value = -70089

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6300.769
$ws.Range("I61").Value = 4531.3
$ws.Range("J61").Value = 12199
$ws.Range("K61").Value = 4531.3
$ws.Range("L61").Value = 12199
$ws.Range("M61").Value = -4329.3
$ws.Range("N61").Value = -12603
$ws.Range("H113").Value = 6300.769
$ws.Range("I113").Value = 4531.3
$ws.Range("J113").Value = 12199
$ws.Range("K113").Value = 4531.3
$ws.Range("L113").Value = 12199
$ws.Range("M113").Value = -2361.3
$ws.Range("N113").Value = -16539
$ws.Range("H132").Value = 7667.778
$ws.Range("I132").Value = 5000
$ws.Range("J132").Value = 9001.666999999999
$ws.Range("K132").Value = 15000
$ws.Range("L132").Value = 27005.001
$ws.Range("M132").Value = -12470
$ws.Range("N132").Value = -32065.001
$ws.Range("H136").Value = 5328.5415
$ws.Range("I136").Value = 3363.2354
$ws.Range("K136").Value = 10089.7062
$ws.Range("M136").Value = -7539.706200000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 140000
$ws.Range("I2").Value = 140000
$ws.Range("K2").Value = 140000
$ws.Range("M2").Value = -139888
